$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 24,13
$arr[0,0] = 35.91501269785719
$arr[0,1] = 21.05472306107833
$arr[0,2] = 11.81038654205563
$arr[0,3] = 10.25996335714982
$arr[0,4] = 0
$arr[0,5] = 3.857632238152332
$arr[0,6] = 0
$arr[0,7] = 46.05477534902439
$arr[0,8] = 6.684132768814813
$arr[0,9] = 0
$arr[0,10] = 16.42363388547086
$arr[0,11] = 0
$arr[0,12] = 24.34213714278981
$arr[1,0] = 35.67391653507793
$arr[1,1] = 20.74544779757803
$arr[1,2] = 11.82855766992023
$arr[1,3] = 10.27656689076776
$arr[1,4] = 0
$arr[1,5] = 3.863338655376271
$arr[1,6] = 0
$arr[1,7] = 46.00876633208192
$arr[1,8] = 6.665841935186831
$arr[1,9] = 0
$arr[1,10] = 16.4322787666449
$arr[1,11] = 0
$arr[1,12] = 24.38204253554182
$arr[2,0] = 35.5357440405914
$arr[2,1] = 20.5604314435417
$arr[2,2] = 11.84170954801274
$arr[2,3] = 10.28737201605524
$arr[2,4] = 0
$arr[2,5] = 3.867018243009843
$arr[2,6] = 0
$arr[2,7] = 45.98971517963759
$arr[2,8] = 6.654353698137309
$arr[2,9] = 0
$arr[2,10] = 16.44085626744596
$arr[2,11] = 0
$arr[2,12] = 24.40849891579232
$arr[3,0] = 35.48195996014999
$arr[3,1] = 20.48635392657211
$arr[3,2] = 11.84756987604977
$arr[3,3] = 10.29192919367982
$arr[3,4] = 0
$arr[3,5] = 3.86856212336413
$arr[3,6] = 0
$arr[3,7] = 45.98426012058769
$arr[3,8] = 6.649607146841927
$arr[3,9] = 0
$arr[3,10] = 16.44517225111263
$arr[3,11] = 0
$arr[3,12] = 24.41977051677934
$arr[4,0] = 35.47318264966792
$arr[4,1] = 20.47413562473184
$arr[4,2] = 11.84857320452143
$arr[4,3] = 10.29269522525546
$arr[4,4] = 0
$arr[4,5] = 3.868821172321689
$arr[4,6] = 0
$arr[4,7] = 45.98349357355515
$arr[4,8] = 6.648815052951812
$arr[4,9] = 0
$arr[4,10] = 16.44593842698471
$arr[4,11] = 0
$arr[4,12] = 24.42167174778676
$arr[5,0] = 35.53500842397684
$arr[5,1] = 20.55942695301643
$arr[5,2] = 11.84178655563693
$arr[5,3] = 10.28743285154465
$arr[5,4] = 0
$arr[5,5] = 3.867038884197037
$arr[5,6] = 0
$arr[5,7] = 45.98963227108643
$arr[5,8] = 6.654289947982401
$arr[5,9] = 0
$arr[5,10] = 16.44091115404662
$arr[5,11] = 0
$arr[5,12] = 24.40864894419513
$arr[6,0] = 35.82986291422418
$arr[6,1] = 20.94712762484483
$arr[6,2] = 11.81623749947967
$arr[6,3] = 10.26556184280313
$arr[6,4] = 0
$arr[6,5] = 3.859563442895485
$arr[6,6] = 0
$arr[6,7] = 46.03699989670122
$arr[6,8] = 6.677879412829624
$arr[6,9] = 0
$arr[6,10] = 16.42593519073921
$arr[6,11] = 0
$arr[6,12] = 24.35549028980319
$arr[7,0] = 36.48423574539426
$arr[7,1] = 21.74200461079514
$arr[7,2] = 11.78199592939938
$arr[7,3] = 10.22749477009421
$arr[7,4] = 0
$arr[7,5] = 3.846289801909534
$arr[7,6] = 0
$arr[7,7] = 46.20302855065174
$arr[7,8] = 6.722105944895818
$arr[7,9] = 0
$arr[7,10] = 16.42257013858458
$arr[7,11] = 0
$arr[7,12] = 24.2667992212551
$arr[8,0] = 37.00836729418319
$arr[8,1] = 22.34147508576067
$arr[8,2] = 11.76655028882767
$arr[8,3] = 10.20243598856209
$arr[8,4] = 0
$arr[8,5] = 3.837369210233996
$arr[8,6] = 0
$arr[8,7] = 46.36973498490381
$arr[8,8] = 6.753378033529753
$arr[8,9] = 0
$arr[8,10] = 16.43601940605278
$arr[8,11] = 0
$arr[8,12] = 24.21118435527972
$arr[9,0] = 37.25549673362559
$arr[9,1] = 22.61635460565452
$arr[9,2] = 11.76164034673846
$arr[9,3] = 10.19166133064734
$arr[9,4] = 0
$arr[9,5] = 3.833488726345462
$arr[9,6] = 0
$arr[9,7] = 46.45528538759513
$arr[9,8] = 6.76734324200323
$arr[9,9] = 0
$arr[9,10] = 16.44560402606224
$arr[9,11] = 0
$arr[9,12] = 24.18797125869058
$arr[10,0] = 37.35026352782646
$arr[10,1] = 22.72065912716846
$arr[10,2] = 11.76008589711184
$arr[10,3] = 10.18767058378277
$arr[10,4] = 0
$arr[10,5] = 3.83204459596665
$arr[10,6] = 0
$arr[10,7] = 46.48907516604144
$arr[10,8] = 6.772594282971401
$arr[10,9] = 0
$arr[10,10] = 16.44973207205572
$arr[10,11] = 0
$arr[10,12] = 24.17948237496434
$arr[11,0] = 37.3298022202136
$arr[11,1] = 22.69818734126365
$arr[11,2] = 11.76040711059719
$arr[11,3] = 10.18852609435026
$arr[11,4] = 0
$arr[11,5] = 3.832354492230624
$arr[11,6] = 0
$arr[11,7] = 46.48173601547916
$arr[11,8] = 6.771465032283574
$arr[11,9] = 0
$arr[11,10] = 16.44882085202581
$arr[11,11] = 0
$arr[11,12] = 24.18129717996948
$arr[12,0] = 37.26326991123002
$arr[12,1] = 22.62493208200134
$arr[12,2] = 11.76150634917033
$arr[12,3] = 10.19133122071098
$arr[12,4] = 0
$arr[12,5] = 3.833369410503339
$arr[12,6] = 0
$arr[12,7] = 46.45803738531892
$arr[12,8] = 6.767775994333662
$arr[12,9] = 0
$arr[12,10] = 16.44593365169356
$arr[12,11] = 0
$arr[12,12] = 24.18726682428051
$arr[13,0] = 37.22266915409218
$arr[13,1] = 22.58008602512001
$arr[13,2] = 11.76221937732318
$arr[13,3] = 10.19306106852305
$arr[13,4] = 0
$arr[13,5] = 3.833994369172411
$arr[13,6] = 0
$arr[13,7] = 46.44370266786381
$arr[13,8] = 6.765511502889463
$arr[13,9] = 0
$arr[13,10] = 16.44423008217836
$arr[13,11] = 0
$arr[13,12] = 24.19096269844512
$arr[14,0] = 36.99238620866214
$arr[14,1] = 22.32354632288603
$arr[14,2] = 11.76691379694924
$arr[14,3] = 10.20315266996039
$arr[14,4] = 0
$arr[14,5] = 3.837626363977189
$arr[14,6] = 0
$arr[14,7] = 46.36433930526083
$arr[14,8] = 6.752460146774516
$arr[14,9] = 0
$arr[14,10] = 16.4354627973653
$arr[14,11] = 0
$arr[14,12] = 24.21274348736426
$arr[15,0] = 36.85329650091154
$arr[15,1] = 22.16665375550928
$arr[15,2] = 11.77033610857333
$arr[15,3] = 10.20950322180014
$arr[15,4] = 0
$arr[15,5] = 3.839899801078811
$arr[15,6] = 0
$arr[15,7] = 46.31813878778836
$arr[15,8] = 6.744386951371406
$arr[15,9] = 0
$arr[15,10] = 16.43097236051115
$arr[15,11] = 0
$arr[15,12] = 24.22664064876396
$arr[16,0] = 36.77411876730594
$arr[16,1] = 22.07662408798174
$arr[16,2] = 11.77250370784539
$arr[16,3] = 10.21321471835117
$arr[16,4] = 0
$arr[16,5] = 3.841224144977589
$arr[16,6] = 0
$arr[16,7] = 46.29248016802153
$arr[16,8] = 6.739718920292042
$arr[16,9] = 0
$arr[16,10] = 16.4287158768942
$arr[16,11] = 0
$arr[16,12] = 24.23483022354695
$arr[17,0] = 36.74745388926959
$arr[17,1] = 22.04618085887602
$arr[17,2] = 11.77327181158328
$arr[17,3] = 10.21448148468984
$arr[17,4] = 0
$arr[17,5] = 3.841675423497732
$arr[17,6] = 0
$arr[17,7] = 46.28394977544475
$arr[17,8] = 6.738134183962691
$arr[17,9] = 0
$arr[17,10] = 16.4280079043626
$arr[17,11] = 0
$arr[17,12] = 24.23763674477119
$arr[18,0] = 36.86801813749467
$arr[18,1] = 22.1833341491466
$arr[18,2] = 11.76995117974907
$arr[18,3] = 10.20882110960747
$arr[18,4] = 0
$arr[18,5] = 3.839656060535841
$arr[18,6] = 0
$arr[18,7] = 46.32296226326257
$arr[18,8] = 6.745248895442357
$arr[18,9] = 0
$arr[18,10] = 16.43141660309265
$arr[18,11] = 0
$arr[18,12] = 24.22514094644653
$arr[19,0] = 37.28278046630648
$arr[19,1] = 22.64644391151966
$arr[19,2] = 11.76117519925072
$arr[19,3] = 10.19050486486036
$arr[19,4] = 0
$arr[19,5] = 3.833070618736864
$arr[19,6] = 0
$arr[19,7] = 46.46496045637836
$arr[19,8] = 6.768860564419779
$arr[19,9] = 0
$arr[19,10] = 16.44676816259936
$arr[19,11] = 0
$arr[19,12] = 24.18550520351687
$arr[20,0] = 37.56072181477339
$arr[20,1] = 22.95031362523637
$arr[20,2] = 11.7572166031452
$arr[20,3] = 10.17905492214125
$arr[20,4] = 0
$arr[20,5] = 3.82891416778088
$arr[20,6] = 0
$arr[20,7] = 46.56588545666956
$arr[20,8] = 6.784074947299635
$arr[20,9] = 0
$arr[20,10] = 16.45970708843205
$arr[20,11] = 0
$arr[20,12] = 24.16135848602793
$arr[21,0] = 37.41177327360928
$arr[21,1] = 22.78805501955755
$arr[21,2] = 11.75916663566725
$arr[21,3] = 10.18511846813027
$arr[21,4] = 0
$arr[21,5] = 3.831119114034911
$arr[21,6] = 0
$arr[21,7] = 46.51127828062978
$arr[21,8] = 6.775974554533112
$arr[21,9] = 0
$arr[21,10] = 16.45253551446789
$arr[21,11] = 0
$arr[21,12] = 24.17408475999057
$arr[22,0] = 36.86136002988139
$arr[22,1] = 22.17579240622461
$arr[22,2] = 11.77012458294628
$arr[22,3] = 10.20912930411486
$arr[22,4] = 0
$arr[22,5] = 3.839766201689229
$arr[22,6] = 0
$arr[22,7] = 46.32077875814375
$arr[22,8] = 6.74485929353143
$arr[22,9] = 0
$arr[22,10] = 16.43121474832915
$arr[22,11] = 0
$arr[22,12] = 24.22581833928936
$arr[23,0] = 36.29935841501833
$arr[23,1] = 21.5238433423931
$arr[23,2] = 11.78955653536539
$arr[23,3] = 10.23727986686331
$arr[23,4] = 0
$arr[23,5] = 3.849733690257686
$arr[23,6] = 0
$arr[23,7] = 46.15026320640271
$arr[23,8] = 6.710357217653325
$arr[23,9] = 0
$arr[23,10] = 16.42068832258454
$arr[23,11] = 0
$arr[23,12] = 24.28912086315821
$range = $ws.Range("B2:N25")
$range.Value = $arr
